# This script applies a re-sort / value refresh of the weekly "Breva" price
# rows (r2..r13) so each row now carries the figures described by the
# upstream diff (dates, volumes, prices, unit + origin labels).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44204
$ws.Range("M2").Value = 110
$ws.Range("N2").Value = 7000
$ws.Range("O2").Value = 7500
$ws.Range("P2").Value = 7318
$ws.Range('Q2').Value = '$/bandeja 7 kilos'
$ws.Range('R2').Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S2").Value = 1045
$ws.Range("T2").Value = 7

# Row 3
$ws.Range("D3").Value = 44550
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 24000
$ws.Range("O3").Value = 24000
$ws.Range("P3").Value = 24000
$ws.Range('Q3').Value = '$/bandeja 7 kilos'
$ws.Range('R3').Value = 'Región Metropolitana'
$ws.Range("S3").Value = 3429
$ws.Range("T3").Value = 7

# Row 4
$ws.Range("D4").Value = 44553
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 22000
$ws.Range("O4").Value = 22000
$ws.Range("P4").Value = 22000
$ws.Range('Q4').Value = '$/bandeja 6 kilos'
$ws.Range("S4").Value = 3667
$ws.Range("T4").Value = 6

# Row 5
$ws.Range("D5").Value = 44553
$ws.Range("M5").Value = 150
$ws.Range("N5").Value = 18000
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 18000
$ws.Range('Q5').Value = '$/bandeja 6 kilos'
$ws.Range("S5").Value = 3000
$ws.Range("T5").Value = 6

# Row 6
$ws.Range("D6").Value = 44189
$ws.Range('L6').Value = 'Especial'
$ws.Range("M6").Value = 20
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 15000
$ws.Range('R6').Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S6").Value = 2143

# Row 7
$ws.Range("D7").Value = 44189
$ws.Range('L7').Value = 'Primera'
$ws.Range("M7").Value = 30
$ws.Range("N7").Value = 13000
$ws.Range("O7").Value = 13000
$ws.Range("P7").Value = 13000
$ws.Range('Q7').Value = '$/bandeja 7 kilos'
$ws.Range("S7").Value = 1857
$ws.Range("T7").Value = 7

# Row 8
$ws.Range('L8').Value = 'Especial'
$ws.Range("M8").Value = 20
$ws.Range("N8").Value = 22000
$ws.Range("O8").Value = 22000
$ws.Range("P8").Value = 22000
$ws.Range("S8").Value = 3667

# Row 9
$ws.Range("D9").Value = 44558
$ws.Range("M9").Value = 25
$ws.Range("N9").Value = 18000
$ws.Range("O9").Value = 18000
$ws.Range("P9").Value = 18000
$ws.Range('Q9').Value = '$/bandeja 6 kilos'
$ws.Range("S9").Value = 3000
$ws.Range("T9").Value = 6

# Row 12
$ws.Range("D12").Value = 44561
$ws.Range('L12').Value = 'Primera'
$ws.Range("N12").Value = 18000
$ws.Range("O12").Value = 18000
$ws.Range("P12").Value = 18000
$ws.Range("S12").Value = 3000

# Row 13
$ws.Range("D13").Value = 44572
$ws.Range("M13").Value = 65
$ws.Range("N13").Value = 20000
$ws.Range("O13").Value = 20000
$ws.Range("P13").Value = 20000
$ws.Range('R13').Value = 'Región Metropolitana'
$ws.Range("S13").Value = 3333
